$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row for "SIRLY DEL CARMEN FRANCO BERTEL" period 2509 ---
# Push everything from row 18 down by inserting a blank row at 18 (row 17's
# original content - YUREIDIS / 2508 - stays put for now and gets edited below).
$ws.Rows("18").Insert()

# Copy row17's current formatting (the "last row of table" bottom border) down
# into the newly inserted row 18, since that row will become the new last data row.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# Row 17 is no longer the last row of the table, so give it the "middle row"
# formatting that row 16 uses.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# --- Row 17 becomes: CC | 45565246 | SIRLY DEL CARMEN FRANCO BERTEL | 2509 | 56940 | 1423500 ---
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45565246"
$ws.Range("D17").Value = "SIRLY DEL CARMEN FRANCO BERTEL"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# --- Row 18 keeps YUREIDIS, but her period moves from 2508 to 2509 ---
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1043295944"
$ws.Range("D18").Value = "YUREIDIS ALEJANDRA CONEO GALVIS"
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# --- Update summary figures ---
# Cant. Periodos: 1 -> 2 (now two distinct periods: 2508, 2509)
$ws.Range("F13").Value = 2

# VALOR MORA total: 113880 -> 170820 (three line items of 56940 each)
$ws.Range("E11").Value = 170820
